$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Change text of "Retângulo 27" (id=28, currently "Equipes") to "Aluno"
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Retângulo 27") {
        $shp.TextFrame.TextRange.Text = "Aluno"
    }
}

# Delete the connector "Conector reto 47" (id=48)
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Conector reto 47") {
        $shp.Delete()
    }
}

# Delete the diamond "Losango 51" (id=52)
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Losango 51") {
        $shp.Delete()
    }
}

# Delete the rectangle "Retângulo 6" (id=7, which had text "Aluno")
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Retângulo 6") {
        $shp.Delete()
    }
}
